$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.791.72"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +2.29%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.452.47"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.83%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "575.12"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.28%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "160.66"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +3.42%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.613"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +8.49%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.452.78"

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.75%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.125"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +2.32%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.453"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +3.87%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.051.35"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.83%  "

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.62%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000192"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.58%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "28.17"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +3.26%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.851.57"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +2.40%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.456.40"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.18%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.48"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +3.07%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.35"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +2.04%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "380.23"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.25%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "8.12"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.01%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.553"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +3.75%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "72.80"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.90%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.997"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.20%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000119"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.73%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.96"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +4.96%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.178"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.32%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.17%  "

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +10.53%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.15"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.60%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.04"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +2.49%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "23.61"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.36%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.30"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +6.50%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.61"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +10.92%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "161.30"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.13%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.92"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +5.17%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0779"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +3.01%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.917.03"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.27%  "

$ws.Range("B40").Value = "Filecoin"
$ws.Range("C40").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.69"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +7.56%  "

$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.76"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +5.51%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "26.68"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.25%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0322"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +2.11%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "42.99"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +2.42%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.778"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +3.39%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "26.00"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +12.21%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.09"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +3.38%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "319.99"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +10.54%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.110"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +4.99%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.881"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +5.18%  "

$ws.Range("B51").Value = "Cosmos"
$ws.Range("C51").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.60"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +3.11%  "
